$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new "2022" column by duplicating column Q's formatting (copy+insert
# keeps every cell's existing style index intact instead of minting new
# ones), then overwrite the inserted cells with the new year's figures.
$ws.Columns("Q:Q").Copy()
$ws.Columns("R:R").Insert(-4161)

$ws.Range("R4").Value = 2022
$ws.Range("R5").Value = 8.6821914120339212
$ws.Range("R6").Value = 12.221423436376707

# Match the recorded view state: active cell moved one column past the new data.
$ws.Range("S4").Select()
